$p = $ppt.ActivePresentation

# Locate the shape whose text contains the bullet we need to split:
# "...And Overdose Drug" -> "Overdose Drug" + new bullet "Marijuana related death for treatment vs Recreation"
$targetShape = $null
$targetSlide = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.TextRange.Text -like "*And Overdose Drug*") {
                $targetShape = $shp
                $targetSlide = $sl
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Find which paragraph (line) holds "And Overdose Drug" and rewrite its text.
$lineCount = $tr.Lines().Count
$paraIndex = -1
for ($i = 1; $i -le $lineCount; $i++) {
    $ln = $tr.Lines($i, 1)
    if ($ln.Text -like "*And Overdose Drug*") {
        $paraIndex = $i
    }
}

$line = $tr.Lines($paraIndex, 1)
$line.Text = "Overdose Drug "

# Re-fetch the (now shortened) paragraph so Start/Length are current.
$tr = $targetShape.TextFrame.TextRange
$line = $tr.Lines($paraIndex, 1)
$afterLineEnd = $line.Start + $line.Length

# Insert a brand-new paragraph right after it, then add the remaining
# wording as three separate runs (mirrors three distinct typed/edited runs).
$anchor = $tr.Characters($afterLineEnd, 0)
$run1 = $anchor.InsertAfter([char]13 + "Marijuana ")

$anchor2 = $tr.Characters($run1.Start + $run1.Length, 0)
$run2 = $anchor2.InsertAfter("related death for ")

$anchor3 = $tr.Characters($run2.Start + $run2.Length, 0)
$run3 = $anchor3.InsertAfter("treatment vs Recreation  ")
